# JS-SPA-Self-Evaluation-Protocol.xlsx edit script
# Commit: "Packet manager changed to NPM; Admin controllers for editing,
#          creating and deleting added;"
#
# Summary of data changes (per the target OOXML diff):
#   - GitHub commits count (C9) filled in: 46
#   - "AngularJS Project Structure" comment (E12) cleared
#   - "Admin Home Screen" comment (E34) cleared, and its old "Yes" value is
#     instead distributed as per-row Yes/No answers across the Admin
#     options block (C35:C50)
#   - Selection/scroll position moved down to the Admin block (E15 / row 31)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- GitHub section -------------------------------------------------------
# "Numbers of Commits in GitHub" now has an actual value.
$ws.Range("C9").Value = 46

# --- Basic Options section -------------------------------------------------
# Clear the leftover comment on the "AngularJS Project Structure" row - the
# grading rubric no longer carries free-text commentary here.
$ws.Range("E12").ClearContents()

# --- Admin Options section --------------------------------------------------
# The single "Yes" comment that used to live on the section header row (E34)
# is removed ...
$ws.Range("E34").ClearContents()

# ... and each Admin feature row below gets its own explicit Yes/No answer
# in column C (previously blank, i.e. ungraded).
$ws.Range("C35").Value = "Yes"
$ws.Range("C36").Value = "Yes"
$ws.Range("C37").Value = "Yes"
$ws.Range("C38").Value = "Yes"
$ws.Range("C39").Value = "Yes"
$ws.Range("C40").Value = "No"
$ws.Range("C41").Value = "No"
$ws.Range("C42").Value = "No"
$ws.Range("C43").Value = "Yes"
$ws.Range("C44").Value = "No"
$ws.Range("C45").Value = "No"
$ws.Range("C46").Value = "No"
$ws.Range("C47").Value = "Yes"
$ws.Range("C48").Value = "No"
$ws.Range("C49").Value = "No"
$ws.Range("C50").Value = "No"

# --- View state -------------------------------------------------------------
# Scroll/select down into the Admin Options block that was just filled in.
$ws.Range("E15").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
